$d = $word.ActiveDocument

# --- Change 1: "${local_do_evento}" -> "${local_do_evento}," -------------
# The match text lives entirely inside one run, so a plain Find/Replace keeps
# that run's own formatting (bold placeholder) untouched.
$ok1 = $d.Content.Find.Execute('${local_do_evento}', $false, $false, $false, $false, $false,
                                $true, 1, $false, '${local_do_evento},', 2)
if (-not $ok1) {
    throw "edit.ps1: could not find '`${local_do_evento}' to append the comma"
}

# --- Change 2: "},", " ", " ", "São Paulo" (4 runs) -> "}, " + " São Paulo" (2 runs)
# Locate the exact span unambiguously via its surrounding context ("hoje},  S").
$anchor = $d.Content
$ok2 = $anchor.Find.Execute('hoje},  S', $false, $false, $false, $false, $false,
                             $true, 1, $false, "", 0)
if (-not $ok2) {
    throw "edit.ps1: could not find the 'hoje},  S' anchor around the São Paulo run"
}

$commaStart = $anchor.Start + 4   # skip "hoje", land right before "},"
$commaEnd   = $commaStart + 2     # end of "},"
$spaceEnd   = $commaEnd + 1       # end of the first (soon to be deleted) space run

# Delete the lone space run sandwiched between "}," and the second space run.
# Because that second space run and "São Paulo" already share identical
# formatting, Word auto-merges them into a single " São Paulo" run once the
# gap run disappears and they become adjacent.
$gap = $d.Range($commaEnd, $spaceEnd)
$gap.Delete()

# Re-insert a trailing space onto the "}," run itself (via InsertAfter on that
# run's own range) so the text reads "}, " while keeping that run's original
# (non-bold) formatting rather than picking up the neighbour's.
$commaRun = $d.Range($commaStart, $commaEnd)
$commaRun.InsertAfter(" ")
